$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.253.59"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "2.645.31"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.59%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -0.97%  "
$ws.Range("E9").Value = "  -2.29%  "
$ws.Range("E10").Value = "  -1.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.29"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("E12").Value = "  -0.88%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.93"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.01%  "
$ws.Range("D14").Value = "3.128.04"
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000187"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.52%  "
$ws.Range("D16").Value = "68.128.52"
$ws.Range("E16").Value = "  -0.33%  "
$ws.Range("D17").Value = "2.665.20"
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "359.44"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.58%  "
$ws.Range("E20").Value = "  -1.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.06%  "
$ws.Range("E23").Value = "  -0.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.49"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.71"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.97%  "
$ws.Range("D27").Value = "2.778.68"
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("E28").Value = "  -2.79%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "559.91"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.09%  "
$ws.Range("E32").Value = "  -2.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.88"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.64"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.59%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  -2.86%  "
$ws.Range("E37").Value = "  -1.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.66"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.27%  "
$ws.Range("E39").Value = "  -1.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.86"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.33"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.60"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.21%  "
$ws.Range("E43").Value = "  -5.79%  "
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "157.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "21.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.69"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.99%  "
$ws.Range("E49").Value = "  -2.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.574"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.615"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.28%  "
